# Excel COM-interop refresh script ("chore: update Sheets via scheduled runner").
#
# The workbook tracks Leve crafting profitability per job (one worksheet per
# crafting class: ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR). Columns H:N are a market-price
# snapshot (currentAveragePrice[ NQ/HQ], LevePrice[ NQ/HQ], LeveProfit[ NQ/HQ])
# recomputed by the scheduled data-refresh job; columns A:G (item/leve metadata)
# are untouched. This script pokes the refreshed values straight into the cells
# that changed, grouped by worksheet then by row.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart
$ws.Range("H19").Value = 514.4828
$ws.Range("J19").Value = 507.11765
$ws.Range("L19").Value = 507.11765
$ws.Range("N19").Value = -857.11765
# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 1146.1538
$ws.Range("I40").Value = 1150
$ws.Range("J40").Value = 1144.4445
$ws.Range("K40").Value = 1150
$ws.Range("L40").Value = 1144.4445
$ws.Range("M40").Value = -975
$ws.Range("N40").Value = -1494.4445
# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 13335.909
$ws.Range("I62").Value = 16646
$ws.Range("J62").Value = 6242.857
$ws.Range("K62").Value = 16646
$ws.Range("L62").Value = 6242.857
$ws.Range("M62").Value = -16022
$ws.Range("N62").Value = -7490.857
# Row 64: Forged from the Void
$ws.Range("H64").Value = 3967
$ws.Range("J64").Value = 4183
$ws.Range("L64").Value = 4183
$ws.Range("N64").Value = -4679
# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 13335.909
$ws.Range("I65").Value = 16646
$ws.Range("J65").Value = 6242.857
$ws.Range("K65").Value = 83230
$ws.Range("L65").Value = 31214.285
$ws.Range("M65").Value = -80110
$ws.Range("N65").Value = -37454.285
# Row 67: Dodging the Draft (L)
$ws.Range("H67").Value = 3967
$ws.Range("J67").Value = 4183
$ws.Range("L67").Value = 4183
$ws.Range("N67").Value = -5899
# Row 116: Growing Up
$ws.Range("H116").Value = 120521.39
$ws.Range("I116").Value = 178282.08
$ws.Range("K116").Value = 178282.08
$ws.Range("M116").Value = -174840.08
# Row 125: Body over Mind
$ws.Range("H125").Value = 490
$ws.Range("I125").Value = 527.4666999999999
$ws.Range("J125").Value = 419.75
$ws.Range("K125").Value = 4747.2003
$ws.Range("L125").Value = 3777.75
$ws.Range("M125").Value = -2287.2003
$ws.Range("N125").Value = -8697.75
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 2868.5645
$ws.Range("I132").Value = 1505.7142
$ws.Range("J132").Value = 8005.4614
$ws.Range("K132").Value = 4517.142599999999
$ws.Range("L132").Value = 24016.3842
$ws.Range("M132").Value = -1987.142599999999
$ws.Range("N132").Value = -29076.3842
# Row 135: For Tired Minds
$ws.Range("H135").Value = 414.83334
$ws.Range("I135").Value = 395.5
$ws.Range("J135").Value = 540.5
$ws.Range("K135").Value = 3559.5
$ws.Range("L135").Value = 4864.5
$ws.Range("M135").Value = -1024.5
$ws.Range("N135").Value = -9934.5
# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 2641.6
$ws.Range("I137").Value = 2768.7036
$ws.Range("J137").Value = 2377.6155
$ws.Range("K137").Value = 8306.110799999999
$ws.Range("L137").Value = 7132.8465
$ws.Range("M137").Value = -5756.110799999999
$ws.Range("N137").Value = -12232.8465
# Row 141: Remedy for Reason
$ws.Range("H141").Value = 1145.4667
$ws.Range("I141").Value = 1101.1628
$ws.Range("J141").Value = 2098
$ws.Range("K141").Value = 3303.4884
$ws.Range("L141").Value = 6294
$ws.Range("M141").Value = 1876.5116
$ws.Range("N141").Value = -16654

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 3251.5
$ws.Range("I2").Value = 3949.8333
$ws.Range("J2").Value = 1156.5
$ws.Range("K2").Value = 3949.8333
$ws.Range("L2").Value = 1156.5
$ws.Range("M2").Value = -3836.8333
$ws.Range("N2").Value = -1382.5
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 2637.31
$ws.Range("I32").Value = 2497.2268
$ws.Range("J32").Value = 7166.6665
$ws.Range("K32").Value = 2497.2268
$ws.Range("L32").Value = 7166.6665
$ws.Range("M32").Value = -2210.2268
$ws.Range("N32").Value = -7740.6665
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 3877.9688
$ws.Range("J61").Value = 5186.6
$ws.Range("L61").Value = 5186.6
$ws.Range("N61").Value = -5610.6
# Row 63: Rivets Run through It
$ws.Range("H63").Value = 6379.619
$ws.Range("J63").Value = 4577.778
$ws.Range("L63").Value = 4577.778
$ws.Range("N63").Value = -5949.778
# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 6379.619
$ws.Range("J66").Value = 4577.778
$ws.Range("L66").Value = 22888.89
$ws.Range("N66").Value = -29752.89
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 3377.2812
$ws.Range("I74").Value = 1713.375
$ws.Range("J74").Value = 5041.1875
$ws.Range("K74").Value = 1713.375
$ws.Range("L74").Value = 5041.1875
$ws.Range("M74").Value = -839.375
$ws.Range("N74").Value = -6789.1875
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 3377.2812
$ws.Range("I77").Value = 1713.375
$ws.Range("J77").Value = 5041.1875
$ws.Range("K77").Value = 8566.875
$ws.Range("L77").Value = 25205.9375
$ws.Range("M77").Value = -4198.875
$ws.Range("N77").Value = -33941.9375
# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 2023.5625
$ws.Range("I110").Value = 1939.75
$ws.Range("K110").Value = 1939.75
$ws.Range("M110").Value = 105.25
# Row 116: No Scope
$ws.Range("H116").Value = 3251.5
$ws.Range("I116").Value = 3949.8333
$ws.Range("J116").Value = 1156.5
$ws.Range("K116").Value = 3949.8333
$ws.Range("L116").Value = 1156.5
$ws.Range("M116").Value = -1655.8333
$ws.Range("N116").Value = -5744.5
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 3877.9688
$ws.Range("J136").Value = 5186.6
$ws.Range("L136").Value = 15559.8
$ws.Range("N136").Value = -20659.8

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 3251.5
$ws.Range("I3").Value = 3949.8333
$ws.Range("J3").Value = 1156.5
$ws.Range("K3").Value = 3949.8333
$ws.Range("L3").Value = 1156.5
$ws.Range("M3").Value = -3835.8333
$ws.Range("N3").Value = -1384.5
# Row 99: Meddle in Metal
$ws.Range("H99").Value = 10655907
$ws.Range("I99").Value = 3852400.5
$ws.Range("J99").Value = 33334264
$ws.Range("K99").Value = 3852400.5
$ws.Range("L99").Value = 33334264
$ws.Range("M99").Value = -3850902.5
$ws.Range("N99").Value = -33337260

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 3342.362
$ws.Range("I31").Value = 2754.1155
$ws.Range("K31").Value = 2754.1155
$ws.Range("M31").Value = -2459.1155
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 3342.362
$ws.Range("I34").Value = 2754.1155
$ws.Range("K34").Value = 2754.1155
$ws.Range("M34").Value = -2552.1155
# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 1687.16
$ws.Range("I132").Value = 1000.46295
$ws.Range("K132").Value = 3001.38885
$ws.Range("M132").Value = -471.3888499999998
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1312.0735
$ws.Range("I134").Value = 982.2549
$ws.Range("J134").Value = 2301.5293
$ws.Range("K134").Value = 2946.7647
$ws.Range("L134").Value = 6904.5879
$ws.Range("M134").Value = -411.7647000000002
$ws.Range("N134").Value = -11974.5879

$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up
$ws.Range("H12").Value = 97.291664
$ws.Range("I12").Value = 23.857143
$ws.Range("J12").Value = 127.52941
$ws.Range("K12").Value = 71.57142899999999
$ws.Range("L12").Value = 382.58823
$ws.Range("M12").Value = 101.428571
$ws.Range("N12").Value = -728.5882300000001
# Row 92: Oh No Udon
$ws.Range("H92").Value = 241
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 241
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 723
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -3219

$ws = $wb.Worksheets.Item("GSM")
# Row 107: Whetstones for the Workers
$ws.Range("H107").Value = 9585.182000000001
$ws.Range("I107").Value = 12803.875
$ws.Range("J107").Value = 1002
$ws.Range("K107").Value = 12803.875
$ws.Range("L107").Value = 1002
$ws.Range("M107").Value = -10883.875
$ws.Range("N107").Value = -4842
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 1650
$ws.Range("I122").Value = 1375
$ws.Range("J122").Value = 1870
$ws.Range("K122").Value = 4125
$ws.Range("L122").Value = 5610
$ws.Range("M122").Value = -1675
$ws.Range("N122").Value = -10510

$ws = $wb.Worksheets.Item("LTW")
# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 2078.2727
$ws.Range("I82").Value = 1994.3103
$ws.Range("J82").Value = 2240.6
$ws.Range("K82").Value = 1994.3103
$ws.Range("L82").Value = 2240.6
$ws.Range("M82").Value = -1633.3103
$ws.Range("N82").Value = -2962.6
# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 2078.2727
$ws.Range("I85").Value = 1994.3103
$ws.Range("J85").Value = 2240.6
$ws.Range("K85").Value = 1994.3103
$ws.Range("L85").Value = 2240.6
$ws.Range("M85").Value = -746.3103000000001
$ws.Range("N85").Value = -4736.6

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1668.4203
$ws.Range("I132").Value = 876.30615
$ws.Range("J132").Value = 3609.1
$ws.Range("K132").Value = 2628.91845
$ws.Range("L132").Value = 10827.3
$ws.Range("M132").Value = -98.91845000000012
$ws.Range("N132").Value = -15887.3
# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 1786.8572
$ws.Range("I136").Value = 995.4
$ws.Range("J136").Value = 2506.3635
$ws.Range("K136").Value = 2986.2
$ws.Range("L136").Value = 7519.0905
$ws.Range("M136").Value = -436.1999999999998
$ws.Range("N136").Value = -12619.0905
